# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
# (GitHub Actions scheduled update) described by the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.870.23"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").Value = "1.648.01"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  +0.46%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.02"
$ws.Range("E5").Value = "  +0.81%  "
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("E7").Value = "  +0.55%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.252"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.25"
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("D12").Value = "1.647.73"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.17"
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.85"
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("D16").Value = "26.878.23"
$ws.Range("E16").Value = "  +0.69%  "
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "214.82"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.40"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("E21").Value = "  +11.15%  "
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.36"
$ws.Range("E23").Value = "  -1.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "147.50"
$ws.Range("E24").Value = "  +1.31%  "
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.118"
$ws.Range("E26").Value = "  -1.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.21"
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.69"
$ws.Range("E28").Value = "  -0.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0507"
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.18"
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("D33").Value = "1.297.61"
$ws.Range("E33").Value = "  +1.52%  "
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("E35").Value = "  +1.50%  "
$ws.Range("E36").Value = "  -1.84%  "
$ws.Range("E37").Value = "  +0.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.823"
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("E39").Value = "  +0.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.809"
$ws.Range("E40").Value = "  -0.93%  "
$ws.Range("E41").Value = "  -0.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.34"
$ws.Range("E42").Value = "  -2.27%  "
$ws.Range("D43").Value = "1.788.30"
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "61.85"
$ws.Range("E44").Value = "  +3.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.08"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("E46").Value = "  +1.37%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0521"
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.67"
$ws.Range("E48").Value = "  -1.30%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0971"
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.01"
$ws.Range("E50").Value = "  +0.59%  "
$ws.Range("E51").Value = "  +0.21%  "
